$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title
Replace-Text "Unveiling the Mysteries of Quantum Entanglement" "Chemistry: The World Around Us Unveiled"

# Author name (merges "Ethan J" + "." + " Smith" into a single run)
Replace-Text "Ethan J. Smith" "Eleanor Spencer"

# Email user part (single run, no merge)
Replace-Text "ethan" "eleanorspencer304@protonmail"

# Email domain part (merges "smith@universalresearch" + "." + "org" into a single run)
Replace-Text "smith@universalresearch.org" "com"

# Body paragraph, first run-group (before the first line break)
Replace-Text "In the realm of quantum physics, lies a phenomenon that has captivated the minds of scientists and philosophers alike - quantum entanglement. This enigmatic phenomenon arises when two particles become correlated in such a way that their states become interconnected, transcending the boundaries of space and time. Entangled particles exhibit a remarkable property: measuring the state of one particle instantaneously affects the state of the other, regardless of the distance between them, creating a non-local connection that defies classical intuition. Delving into the depths of quantum entanglement unveils a rich tapestry of implications for our understanding of the universe, challenging our notions of locality and revealing the profound interconnectedness of all things." "Chemistry, the science of change and transformation, explores the fundamental principles that govern the composition, structure, and behavior of matter. From the intricate processes occurring within our cells to the grand marvels of the cosmos, chemistry weaves its intricate threads through the fabric of our existence. As we embark on this captivating journey into the realm of chemistry, let us unravel the mysteries that lie hidden within the microscopic universe of atoms and molecules, unveiling the secrets that determine the everyday phenomena that surround us."

# Body paragraph, second run-group (between the two line-break pairs)
Replace-Text "The discovery of quantum entanglement has spurred a revolution in scientific thought, blurring the lines between the separate and the shared, and ushering in a new era of quantum technology. From quantum computing and cryptography to teleportation and more, the potential applications of quantum entanglement are vast and hold the promise of transformative advancements across numerous fields. Yet, despite the remarkable progress made in studying this phenomenon, the underlying mechanisms that govern quantum entanglement remain shrouded in mystery, beckoning us to explore the deepest corners of the quantum realm and unlock the secrets it holds." "In this exploration, we shall delve into the captivating world of chemical reactions, where the rearrangement of atoms and molecules gives rise to new substances with unique properties. We shall unravel the intricate interplay of energy and matter, witnessing the transformation of substances from one state to another. Moreover, we shall uncover the secrets of chemical bonding, the fundamental force that holds atoms together, enabling the formation of diverse and complex molecules."

# Body paragraph, third run-group (after the last line-break pair)
Replace-Text "Finally, the profound implications of quantum entanglement extend beyond the realm of science, touching upon the very core of our existence. It invites us to contemplate the nature of reality, the interconnectedness of all things, and the very fabric of spacetime. Quantum entanglement raises philosophical questions that challenge our understanding of the universe and our place within it, inspiring new perspectives on consciousness, free will, and the fundamental unity of all life. The study of quantum entanglement thus becomes a profound journey of discovery, not only into the intricacies of the physical world but also into the depths of our own consciousness and the nature of reality itself." "Furthermore, we shall explore the practical applications of chemistry in our daily lives, from the production of medicines and materials to the purification of water and the generation of energy. We shall witness how chemistry touches every aspect of our existence, from the food we eat and the clothes we wear to the technologies that shape our modern world."

# Summary heading (drop the stray lastRenderedPageBreak marker)
Replace-Text "Summary" "Summary"

# Summary paragraph (single run-group, all three sentences)
Replace-Text "Quantum entanglement, a mysterious phenomenon in quantum physics, reveals the interconnectedness of particles beyond the limits of space and time. It challenges classical notions of locality and has profound implications for our comprehension of the universe. The exploration of quantum entanglement opens doors to revolutionary technological advancements and invites philosophical contemplation of reality, unity, and consciousness, offering a window into the deepest mysteries of the quantum realm and the very essence of existence." "Through this exploration, we have gained a deeper understanding of chemistry, its fundamental principles, and its pervasive influence on our world. From the microscopic interactions of atoms and molecules to the macroscopic phenomena that shape our lives, chemistry reveals the intricate symphony of the material world. By unraveling the enigmas of chemistry, we unlock the potential for innovation, progress, and a profound appreciation for the wonders that surround us."

# Add a new empty paragraph at the very end of the document body
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

Write-Host "done"
